$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.336682677268982
$ws.Range("B1").Value = 2.007334232330322
$ws.Range("C1").Value = 2.519518613815308
$ws.Range("D1").Value = 4.356553077697754
$ws.Range("E1").Value = 1.044628024101257
